# TopUp.cs now checks the existing card balance before topping up, so the
# transactions sheet gains two new columns to record the card used for the
# payment and the resulting balance after the top up was applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the card-balance tracking columns (D and E).
$ws.Range("D1").Value = "Card Number"
$ws.Range("E1").Value = "Card Balance After Payment"

# Replace the single sample transaction (row 2) with the current one:
# student 0, Spring 2025, $20 tuition paid, paid with card 95019925,
# leaving a balance of 68 on the card afterwards.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Spring 2025"
$ws.Range("C2").Value = 20
# Leading apostrophe keeps the digit-only card number stored as text
# instead of being auto-converted to a number.
$ws.Range("D2").Value = "'95019925"
$ws.Range("E2").Value = 68

# Drop the older sample transactions that used to live in rows 3-4.
$ws.Range("A3:E4").Clear()

$ws.Range("E1").Select()
